$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.870.84"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "3.104.81"
$ws.Range("E3").Value = "  +1.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "3.101.04"
$ws.Range("E8").Value = "  +1.05%  "

# Row 9
$ws.Range("E9").Value = "  +1.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "

# Row 11
$ws.Range("E11").Value = "  +1.22%  "

# Row 12
$ws.Range("E12").Value = "  +4.19%  "

# Row 13
$ws.Range("D13").Value = "3.637.57"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14
$ws.Range("E14").Value = "  +2.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.92%  "

# Row 17
$ws.Range("D17").Value = "57.909.76"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$ws.Range("D18").Value = "3.097.07"
$ws.Range("E18").Value = "  +1.65%  "

# Row 19
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "341.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.77%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.512"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.55%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "

# Row 27
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0916"
$ws.Range("E28").Value = "  +1.35%  "

# Row 29
$ws.Range("E29").Value = "  +0.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("E31").Value = "  +1.48%  "

# Row 32
$ws.Range("E32").Value = "  +4.19%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.26%  "

# Row 39
$ws.Range("E39").Value = "  -0.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0667"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

# Row 41
$ws.Range("E41").Value = "  +14.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.683"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.40%  "

# Row 44
$ws.Range("D44").Value = "3.144.77"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("D47").Value = "2.286.32"
$ws.Range("E47").Value = "  +0.49%  "

# Row 48
$ws.Range("E48").Value = "  +3.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "

# Row 51
$ws.Range("E51").Value = "  +2.70%  "
